$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3625657.5
$ws.Range("I18").Value = 876.5
$ws.Range("K18").Value = 876.5
$ws.Range("M18").Value = -592.5
$ws.Range("H32").Value = 8997.666999999999
$ws.Range("J32").Value = 8997.666999999999
$ws.Range("L32").Value = 8997.666999999999
$ws.Range("N32").Value = -9649.666999999999
$ws.Range("H55").Value = 196
$ws.Range("I55").Value = 156
$ws.Range("J55").Value = 284
$ws.Range("K55").Value = 156
$ws.Range("L55").Value = 284
$ws.Range("M55").Value = 58
$ws.Range("N55").Value = -712
$ws.Range("H100").Value = 851.6923
$ws.Range("I100").Value = 864.3333
$ws.Range("J100").Value = 700
$ws.Range("K100").Value = 864.3333
$ws.Range("L100").Value = 700
$ws.Range("M100").Value = -323.3333
$ws.Range("N100").Value = -1782
$ws.Range("H116").Value = 7847
$ws.Range("J116").Value = 7988.5
$ws.Range("L116").Value = 7988.5
$ws.Range("N116").Value = -14872.5
$ws.Range("H137").Value = 2570.5483
$ws.Range("I137").Value = 2160.4546
$ws.Range("J137").Value = 3573
$ws.Range("K137").Value = 6481.3638
$ws.Range("L137").Value = 10719
$ws.Range("M137").Value = -3931.3638
$ws.Range("N137").Value = -15819
$ws.Range("H141").Value = 1864.4642
$ws.Range("I141").Value = 934.6458
$ws.Range("K141").Value = 2803.9374
$ws.Range("M141").Value = 2376.0626

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9214
$ws.Range("I32").Value = 4596.525
$ws.Range("K32").Value = 4596.525
$ws.Range("M32").Value = -4309.525
$ws.Range("H45").Value = 84617800
$ws.Range("I45").Value = 122224230
$ws.Range("K45").Value = 122224230
$ws.Range("M45").Value = -122223853
$ws.Range("H61").Value = 11549.125
$ws.Range("I61").Value = 9899
$ws.Range("J61").Value = 16499.5
$ws.Range("K61").Value = 9899
$ws.Range("L61").Value = 16499.5
$ws.Range("M61").Value = -9687
$ws.Range("N61").Value = -16923.5
$ws.Range("H63").Value = 3536.4348
$ws.Range("I63").Value = 2544.4
$ws.Range("K63").Value = 2544.4
$ws.Range("M63").Value = -1858.4
$ws.Range("H66").Value = 3536.4348
$ws.Range("I66").Value = 2544.4
$ws.Range("K66").Value = 12722
$ws.Range("M66").Value = -9290
$ws.Range("H74").Value = 4519.2915
$ws.Range("I74").Value = 2197
$ws.Range("J74").Value = 9163.875
$ws.Range("K74").Value = 2197
$ws.Range("L74").Value = 9163.875
$ws.Range("M74").Value = -1323
$ws.Range("N74").Value = -10911.875
$ws.Range("H77").Value = 4519.2915
$ws.Range("I77").Value = 2197
$ws.Range("J77").Value = 9163.875
$ws.Range("K77").Value = 10985
$ws.Range("L77").Value = 45819.375
$ws.Range("M77").Value = -6617
$ws.Range("N77").Value = -54555.375
$ws.Range("H136").Value = 11549.125
$ws.Range("I136").Value = 9899
$ws.Range("J136").Value = 16499.5
$ws.Range("K136").Value = 29697
$ws.Range("L136").Value = 49498.5
$ws.Range("M136").Value = -27147
$ws.Range("N136").Value = -54598.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1383145
$ws.Range("I105").Value = 1821727.5
$ws.Range("J105").Value = 4742.857
$ws.Range("K105").Value = 1821727.5
$ws.Range("L105").Value = 4742.857
$ws.Range("M105").Value = -1819980.5
$ws.Range("N105").Value = -8236.857

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2151.975
$ws.Range("J58").Value = 4166
$ws.Range("L58").Value = 4166
$ws.Range("N58").Value = -4572
$ws.Range("H134").Value = 2202.4194
$ws.Range("J134").Value = 3953.8572
$ws.Range("L134").Value = 11861.5716
$ws.Range("N134").Value = -16931.5716
$ws.Range("H136").Value = 2151.975
$ws.Range("J136").Value = 4166
$ws.Range("L136").Value = 12498
$ws.Range("N136").Value = -17598

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 897.3333
$ws.Range("J2").Value = 1997.5
$ws.Range("L2").Value = 11985
$ws.Range("N2").Value = -12211
$ws.Range("H4").Value = 4558852.5
$ws.Range("I4").Value = 4707480.5
$ws.Range("J4").Value = 100002
$ws.Range("K4").Value = 14122441.5
$ws.Range("L4").Value = 300006
$ws.Range("M4").Value = -14122329.5
$ws.Range("N4").Value = -300230
$ws.Range("H5").Value = 2207.3
$ws.Range("I5").Value = 2164.2222
$ws.Range("J5").Value = 2242.5454
$ws.Range("K5").Value = 6492.6666
$ws.Range("L5").Value = 6727.6362
$ws.Range("M5").Value = -6380.6666
$ws.Range("N5").Value = -6951.6362
$ws.Range("H11").Value = 73007.14
$ws.Range("I11").Value = 92254.55
$ws.Range("K11").Value = 276763.65
$ws.Range("M11").Value = -276623.65
$ws.Range("H13").Value = 99
$ws.Range("I13").Value = 99
$ws.Range("K13").Value = 297
$ws.Range("M13").Value = -129
$ws.Range("H20").Value = 2500
$ws.Range("J20").Value = 2500
$ws.Range("L20").Value = 7500
$ws.Range("N20").Value = -7954
$ws.Range("H26").Value = 272.6
$ws.Range("I26").Value = 278.66666
$ws.Range("J26").Value = 218
$ws.Range("K26").Value = 835.9999799999999
$ws.Range("L26").Value = 654
$ws.Range("M26").Value = -547.9999799999999
$ws.Range("N26").Value = -1230
$ws.Range("H29").Value = 310
$ws.Range("I29").Value = 319.33334
$ws.Range("J29").Value = 282
$ws.Range("K29").Value = 958.0000200000001
$ws.Range("L29").Value = 846
$ws.Range("M29").Value = -681.0000200000001
$ws.Range("N29").Value = -1400
$ws.Range("H40").Value = 235.81818
$ws.Range("I40").Value = 199.22223
$ws.Range("J40").Value = 400.5
$ws.Range("K40").Value = 796.88892
$ws.Range("L40").Value = 1602
$ws.Range("M40").Value = -727.88892
$ws.Range("N40").Value = -1740
$ws.Range("H99").Value = 2283.3333
$ws.Range("I99").Value = 925
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 2775
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = -529
$ws.Range("N99").Value = -19492
$ws.Range("H132").Value = 2326.5
$ws.Range("I132").Value = 3153
$ws.Range("K132").Value = 28377
$ws.Range("M132").Value = -25847
$ws.Range("H135").Value = 2207.3
$ws.Range("I135").Value = 2164.2222
$ws.Range("J135").Value = 2242.5454
$ws.Range("K135").Value = 19477.9998
$ws.Range("L135").Value = 20182.9086
$ws.Range("M135").Value = -16942.9998
$ws.Range("N135").Value = -25252.9086

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2089
$ws.Range("I80").Value = 3050
$ws.Range("J80").Value = 1448.3334
$ws.Range("K80").Value = 3050
$ws.Range("L80").Value = 1448.3334
$ws.Range("M80").Value = -2052
$ws.Range("N80").Value = -3444.3334
$ws.Range("H83").Value = 2089
$ws.Range("I83").Value = 3050
$ws.Range("J83").Value = 1448.3334
$ws.Range("K83").Value = 15250
$ws.Range("L83").Value = 7241.666999999999
$ws.Range("M83").Value = -10258
$ws.Range("N83").Value = -17225.667
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H122").Value = 6504.3335
$ws.Range("I122").Value = 4486.2173
$ws.Range("J122").Value = 9405.375
$ws.Range("K122").Value = 13458.6519
$ws.Range("L122").Value = 28216.125
$ws.Range("M122").Value = -11008.6519
$ws.Range("N122").Value = -33116.125
$ws.Range("H126").Value = 3448.842
$ws.Range("I126").Value = 3281.2354
$ws.Range("K126").Value = 9843.706200000001
$ws.Range("M126").Value = -7373.706200000001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8257.799999999999
$ws.Range("I7").Value = 13444.5
$ws.Range("K7").Value = 13444.5
$ws.Range("M7").Value = -13332.5
$ws.Range("H22").Value = 2732.7778
$ws.Range("J22").Value = 3367.5
$ws.Range("L22").Value = 3367.5
$ws.Range("N22").Value = -3957.5
$ws.Range("H27").Value = 2732.7778
$ws.Range("J27").Value = 3367.5
$ws.Range("L27").Value = 3367.5
$ws.Range("N27").Value = -3581.5
$ws.Range("H40").Value = 5996.1333
$ws.Range("I40").Value = 5495.1665
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 5495.1665
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -5359.1665
$ws.Range("N40").Value = -8272
$ws.Range("H126").Value = 8257.799999999999
$ws.Range("I126").Value = 13444.5
$ws.Range("K126").Value = 40333.5
$ws.Range("M126").Value = -37863.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1040.5312
$ws.Range("J107").Value = 1082.6666
$ws.Range("L107").Value = 3247.9998
$ws.Range("N107").Value = -7087.9998
